$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported between the existing rows 90 and
# 91 (old numbering). Insert a fresh row at position 91 which pushes the
# old rows 91..117 down to 92..118, then populate the new row with its
# data.
$ws.Rows.Item(91).Insert()

$ws.Cells.Item(91, 1).Value = 6
$ws.Cells.Item(91, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(91, 3).Value = "Metropolitana"
$ws.Cells.Item(91, 4).Value = 44463
$ws.Cells.Item(91, 5).Value = 13
$ws.Cells.Item(91, 6).Value = 100112026
$ws.Cells.Item(91, 7).Value = "Haba"
$ws.Cells.Item(91, 8).Value = "Sin especificar"
$ws.Cells.Item(91, 9).Value = "Primera"
$ws.Cells.Item(91, 10).Value = 870
$ws.Cells.Item(91, 11).Value = 8000
$ws.Cells.Item(91, 12).Value = 9000
$ws.Cells.Item(91, 13).Value = 8517
$ws.Cells.Item(91, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(91, 15).Value = "Región Metropolitana"
$ws.Cells.Item(91, 16).Value = 341
$ws.Cells.Item(91, 17).Value = 25
$ws.Cells.Item(91, 18).Value = "Hortaliza"
